# Applies the "Corrected excel sheets for application fix issues" edit.
#
# Summary of changes:
#  - Summary!F2: 963.77 -> 0
#  - Summary!A3: 178.1 -> 211.19
#  - Summary!E3: 81.11 -> 114.2
#  - Repayment schedule!F4: 963.77 -> 921.65
#  - Repayment schedule!G4: 3169.45 -> 3211.57
#  - Repayment schedule!F5: 921.65 -> 932.09
#  - Repayment schedule!G5: 2247.8000000000002 -> 2279.48
#  - Repayment schedule!H5: formula (=G4*(12%/365)*B5) -> static value 31.68
#  - Repayment schedule!F6: 941.6 -> 940.54
#  - Repayment schedule!G6: 1306.2 -> 1338.94
#  - Repayment schedule!H6: 22.17 -> 23.23
#  - Repayment schedule!F7: 950.46 -> 950.56
#  - Repayment schedule!G7: 355.74 -> 388.38
#  - Repayment schedule!H7: 13.31 -> 13.21
#  - Repayment schedule!F8: 355.74 -> 388.38
#  - Repayment schedule!H8: 3.51 -> 3.96
#  - Repayment schedule!K8: 359.25 -> 392.34
#  - Repayment schedule!P8: 359.25 -> 392.34
#  - Repayment schedule column O (O2:O8) data cleared out entirely
#  - Repayment schedule!D8/E8: gain blank formatted cells (matching row's style)
#  - Transactions!A2: 182 -> 59
#  - Transactions!A3: 180 -> 57
#  - Selection/cursor positions updated on Summary, Repayment schedule,
#    Transactions sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate() | Out-Null

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 211.19
$wsSummary.Range("E3").Value = 114.2

$wsSummary.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate() | Out-Null

# Row 4
$wsRepay.Range("F4").Value = 921.65
$wsRepay.Range("G4").Value = 3211.57

# Row 5 - also turn the H5 formula into a plain number, keeping the
# surrounding column's formatting (style index 6) instead of the
# formula cell's own number-format style.
$wsRepay.Range("F5").Value = 932.09
$wsRepay.Range("G5").Value = 2279.48
$wsRepay.Range("H4").Copy() | Out-Null
$wsRepay.Range("H5").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("H5").Value = 31.68

# Row 6
$wsRepay.Range("F6").Value = 940.54
$wsRepay.Range("G6").Value = 1338.94
$wsRepay.Range("H6").Value = 23.23

# Row 7
$wsRepay.Range("F7").Value = 950.56
$wsRepay.Range("G7").Value = 388.38
$wsRepay.Range("H7").Value = 13.21

# Row 8 - gains blank formatted D8/E8 cells (copy formatting from the row
# above, which already carries the right style), and several values change.
$wsRepay.Range("D7").Copy() | Out-Null
$wsRepay.Range("D8").PasteSpecial(-4122) | Out-Null
$wsRepay.Range("E7").Copy() | Out-Null
$wsRepay.Range("E8").PasteSpecial(-4122) | Out-Null

$wsRepay.Range("F8").Value = 388.38
$wsRepay.Range("H8").Value = 3.96
$wsRepay.Range("K8").Value = 392.34
$wsRepay.Range("P8").Value = 392.34

# Column O's data (O2:O8) is removed entirely, leaving just the O1 header.
$wsRepay.Range("O2:O8").Clear() | Out-Null

$wsRepay.Rows.Item(9).Select() | Out-Null

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate() | Out-Null

$wsTrans.Range("A2").Value = 59
$wsTrans.Range("A3").Value = 57

$wsTrans.Range("A2:L3").Select() | Out-Null
